$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep a handle on the existing "Hyperlink" cell style (cellXf index 1) so we
# can re-apply it later without the engine minting a brand-new (duplicate)
# cellXf entry the way Hyperlinks.Add does.
$hyperlinkStyle = $ws.Range("B1").Style

# Update the participant's name in place (same shared-string slot gets rewritten).
$ws.Range("A2").Value = "Mohammad Ezzeddin Pratama"

# Move the little Nama/Email table from A1:B3 down to C12:D14.
$ws.Range("A1:B3").Cut($ws.Range("C12"))

# The cut above leaves the two old hyperlinks still pointing at the
# now-blank B2/B3 cells - drop them before re-creating them on the new
# destination cells.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("B3").Hyperlinks.Delete()

# Recreate the hyperlinks on the moved cells (rebuilds the relationships).
# TextToDisplay here seeds the "display" attribute (with the mailto: prefix,
# matching the source data), but it also clobbers the cell's own text - so
# the plain e-mail address gets written back into the cells afterwards.
$ws.Hyperlinks.Add($ws.Range("D14"), "mailto:billiardo985@gmail.com", "", "", "mailto:billiardo985@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:ezzeddinpratama04@gmail.com", "", "", "mailto:ezzeddinpratama04@gmail.com")

# Restore the plain e-mail text in the cells (Hyperlinks.Add overwrote it
# with the mailto:-prefixed TextToDisplay above).
$ws.Range("D14").Value = "billiardo985@gmail.com"
$ws.Range("D13").Value = "ezzeddinpratama04@gmail.com"

# Hyperlinks.Add reformats the cell itself - restore the original Hyperlink
# cell style so we don't leave the sheet with a second, equivalent style.
$ws.Range("D13").Style = $hyperlinkStyle
$ws.Range("D14").Style = $hyperlinkStyle

# Match the saved selection from the edit.
$ws.Range("C5").Select() | Out-Null
